$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.243.77"
$ws.Range("E2").Value = "'  +0.32%  "

$ws.Range("D3").Value = "'1.857.53"
$ws.Range("E3").Value = "'  +0.28%  "

$ws.Range("E4").Value = "'  +0.09%  "

$ws.Range("D5").Value = "'0.7019"
$ws.Range("E5").Value = "'  +1.91%  "

$ws.Range("D6").Value = "'238.10"
$ws.Range("E6").Value = "'  +0.31%  "

$ws.Range("E7").Value = "'  +0.07%  "

$ws.Range("D8").Value = "'0.08025"
$ws.Range("E8").Value = "'  +3.93%  "

$ws.Range("D9").Value = "'0.3021"
$ws.Range("E9").Value = "'  -0.57%  "

$ws.Range("D10").Value = "'23.48"
$ws.Range("E10").Value = "'  +1.30%  "

$ws.Range("D11").Value = "'0.08184"
$ws.Range("E11").Value = "'  +0.29%  "

$ws.Range("D12").Value = "'1.861.28"
$ws.Range("E12").Value = "'  +0.68%  "

$ws.Range("D13").Value = "'5.197"
$ws.Range("E13").Value = "'  -0.05%  "

$ws.Range("D14").Value = "'0.7069"
$ws.Range("E14").Value = "'  -2.29%  "

$ws.Range("E15").Value = "'  +0.58%  "

$ws.Range("D16").Value = "'29.293.22"
$ws.Range("E16").Value = "'  +0.48%  "

$ws.Range("D17").Value = "'5.825"
$ws.Range("E17").Value = "'  +1.62%  "

$ws.Range("D18").Value = "'0.000007892"
$ws.Range("E18").Value = "'  +1.13%  "

$ws.Range("D19").Value = "'13.27"
$ws.Range("E19").Value = "'  +0.74%  "

$ws.Range("D20").Value = "'237.73"
$ws.Range("E20").Value = "'  +1.44%  "

$ws.Range("B21").Value = "'Dai"
$ws.Range("C21").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'0.9989"
$ws.Range("E21").Value = "'  -0.17%  "

$ws.Range("B22").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "'2.120.46"
$ws.Range("E22").Value = "'  +1.11%  "

$ws.Range("E23").Value = "'  +0.05%  "

$ws.Range("D24").Value = "'7.450"
$ws.Range("E24").Value = "'  -0.75%  "

$ws.Range("D25").Value = "'162.83"
$ws.Range("E25").Value = "'  +0.76%  "

$ws.Range("D26").Value = "'8.888"
$ws.Range("E26").Value = "'  -0.81%  "

$ws.Range("D27").Value = "'0.1430"
$ws.Range("E27").Value = "'  +0.42%  "

$ws.Range("D28").Value = "'18.08"
$ws.Range("E28").Value = "'  +0.03%  "

$ws.Range("D29").Value = "'1.922"
$ws.Range("E29").Value = "'  -2.06%  "

$ws.Range("E30").Value = "'  +1.03%  "

$ws.Range("D31").Value = "'1.474"
$ws.Range("E31").Value = "'  -0.59%  "

$ws.Range("D32").Value = "'4.373"
$ws.Range("E32").Value = "'  -3.19%  "

$ws.Range("D33").Value = "'4.027"
$ws.Range("E33").Value = "'  +0.57%  "

$ws.Range("D34").Value = "'0.05195"
$ws.Range("E34").Value = "'  +0.02%  "

$ws.Range("D35").Value = "'1.162"
$ws.Range("E35").Value = "'  -1.37%  "

$ws.Range("D36").Value = "'0.7176"
$ws.Range("E36").Value = "'  +2.09%  "

$ws.Range("D37").Value = "'1.002"
$ws.Range("E37").Value = "'  -1.78%  "

$ws.Range("D38").Value = "'2.685"
$ws.Range("E38").Value = "'  +1.32%  "

$ws.Range("E39").Value = "'  +0.06%  "

$ws.Range("E40").Value = "'  +1.81%  "

$ws.Range("D41").Value = "'0.9395"
$ws.Range("E41").Value = "'  +3.16%  "

$ws.Range("D42").Value = "'1.152.98"
$ws.Range("E42").Value = "'  +5.37%  "

$ws.Range("D43").Value = "'5.997"
$ws.Range("E43").Value = "'  +0.10%  "

$ws.Range("D44").Value = "'0.4262"
$ws.Range("E44").Value = "'  -0.31%  "

$ws.Range("D45").Value = "'70.57"
$ws.Range("E45").Value = "'  +0.22%  "

$ws.Range("D46").Value = "'1.001"
$ws.Range("E46").Value = "'  +0.11%  "

$ws.Range("D47").Value = "'102.87"
$ws.Range("E47").Value = "'  +0.36%  "

$ws.Range("D48").Value = "'0.5294"
$ws.Range("E48").Value = "'  -3.81%  "

$ws.Range("D49").Value = "'2.050.64"
$ws.Range("E49").Value = "'  +2.70%  "

$ws.Range("D50").Value = "'1.756"
$ws.Range("E50").Value = "'  +0.09%  "

$ws.Range("D51").Value = "'9.144"
$ws.Range("E51").Value = "'  +0.07%  "
